# -----------------------------------------------------------------------
# Edit: add a "Player Info" sheet (new first sheet), add an
# "ODI Batting Extra" sheet (new last sheet), and on the existing
# "ODI Batting" / "ODI Bowling" sheets replace the MATCH_CARD_LINK
# column (a full howstat.com URL) with a MATCH_CODE column holding just
# the numeric match code that used to be the `MatchCode=` query value.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New first sheet: "Player Info"
# ---------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$wsInfo = $wb.Worksheets.Add($firstSheet, $null)
$wsInfo.Name = "Player Info"

$wsInfo.Range("A1").Value = "ID"
$wsInfo.Range("B1").Value = "NAME"
$wsInfo.Range("C1").Value = "BATTING_HAND"
$wsInfo.Range("D1").Value = "BOWL_STYLE"
$infoHeader = $wsInfo.Range("A1:D1")
$infoHeader.Borders.LineStyle = 1
$infoHeader.Font.Bold = $true
$infoHeader.HorizontalAlignment = -4108
$infoHeader.VerticalAlignment = -4160

$wsInfo.Range("A2").NumberFormat = "@"
$wsInfo.Range("A2").Value = "3606"
$wsInfo.Range("B2").Value = "Calum Scott Macleod"
$wsInfo.Range("C2").Value = "Right Handed"
$wsInfo.Range("D2").Value = "Right Arm Medium Fast"

# ---------------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK (col D) -> MATCH_CODE
# ---------------------------------------------------------------------
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBatting.Range("D1").Value = "MATCH_CODE"

$battingCodes = @("2858","2952","2960","2964","3299","3301","3307","3308","3428","3482","3483","3522","3523","3552","3554","3556","3604","3611","3613","3637","3674","3675","3676","3721","3724","3725","3733","3753","3761","3764","3774","3782","3787","3880","3919","3920","3933","3934","3979","3980","4048","4049","4077","4078","4090","4091","4111","4113","4118","4120","4140","4142","4153","4155","4158","4161","4165","4290","4302","4363","4364","4365","4366","4381","4384","4386","4461","4462","4510","4512","4513","4515","4569","4570","4572","4575","4576","4578","4581","4604","4610","4612","4617","4625","4629","4631","4632","4635")

$battingRange = $wsBatting.Range("D2:D89")
$battingRange.NumberFormat = "@"
for ($i = 0; $i -lt $battingCodes.Length; $i++) {
    $row = $i + 2
    $wsBatting.Cells.Item($row, 4).Value = $battingCodes[$i]
}

# ---------------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK (col B) -> MATCH_CODE
# ---------------------------------------------------------------------
$wsBowling = $wb.Worksheets.Item("ODI Bowling")
$wsBowling.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @("2952","2960","2964","3522","3523","3552","3554","3604","3611","3613","3675","3676","3880","3919","3920","3979","4049","4078","4091","4111","4113","4118","4120","4153","4155","4158","4363","4364","4366","4381","4386","4510","4512","4569","4575","4576","4578","4581","4631","4632","4635")

$bowlingRange = $wsBowling.Range("B2:B42")
$bowlingRange.NumberFormat = "@"
for ($i = 0; $i -lt $bowlingCodes.Length; $i++) {
    $row = $i + 2
    $wsBowling.Cells.Item($row, 2).Value = $bowlingCodes[$i]
}

# ---------------------------------------------------------------------
# 4. New last sheet: "ODI Batting Extra"
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsExtra = $wb.Worksheets.Add($null, $lastSheet)
$wsExtra.Name = "ODI Batting Extra"

$wsExtra.Range("A1").Value = "MATCH_CODE"
$wsExtra.Range("B1").Value = "BATTING_POSITION"
$wsExtra.Range("C1").Value = "NUM_4"
$wsExtra.Range("D1").Value = "NUM_6"
$wsExtra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$wsExtra.Range("F1").Value = "MAN_OF_MATCH"
$extraHeader = $wsExtra.Range("A1:F1")
$extraHeader.Borders.LineStyle = 1
$extraHeader.Font.Bold = $true
$extraHeader.HorizontalAlignment = -4108
$extraHeader.VerticalAlignment = -4160

$wsExtra.Range("A2:A21").NumberFormat = "@"
$wsExtra.Range("C2:E21").NumberFormat = "@"

$extraData = @(
    @("4510","","","","","NO"),
    @("4512","3","2","0","9.52%","NO"),
    @("4513","","","","","NO"),
    @("4515","","","","","NO"),
    @("4569","3","3","0","11.27%","NO"),
    @("4570","3","2","0","8.84%","NO"),
    @("4572","3","1","0","2.79%","NO"),
    @("4575","3","4","0","19.42%","NO"),
    @("4576","3","0","0","0.33%","NO"),
    @("4578","3","3","0","9.22%","NO"),
    @("4581","3","3","0","11.70%","NO"),
    @("4604","3","5","0","19.38%","NO"),
    @("4610","3","2","0","9.03%","NO"),
    @("4612","3","7","0","26.94%","NO"),
    @("4617","3","10","3","49.23%","YES"),
    @("4625","3","1","0","2.29%","NO"),
    @("4629","4","7","0","29.01%","YES"),
    @("4631","","","","","NO"),
    @("4632","4","8","0","30.31%","NO"),
    @("4635","","","","","NO")
)

for ($i = 0; $i -lt $extraData.Length; $i++) {
    $row = $i + 2
    $rec = $extraData[$i]

    $wsExtra.Cells.Item($row, 1).Value = $rec[0]

    if ($rec[1] -ne "") {
        $wsExtra.Cells.Item($row, 2).Value = [double]$rec[1]
    }
    if ($rec[2] -ne "") {
        $wsExtra.Cells.Item($row, 3).Value = $rec[2]
    }
    if ($rec[3] -ne "") {
        $wsExtra.Cells.Item($row, 4).Value = $rec[3]
    }
    if ($rec[4] -ne "") {
        $wsExtra.Cells.Item($row, 5).Value = $rec[4]
    }
    $wsExtra.Cells.Item($row, 6).Value = $rec[5]
}
